$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5393.1816
$ws.Range("I43").Value = 5853.375
$ws.Range("J43").Value = 4166
$ws.Range("K43").Value = 5853.375
$ws.Range("L43").Value = 4166
$ws.Range("M43").Value = -5784.375
$ws.Range("N43").Value = -4304
$ws.Range("H53").Value = 7122.4116
$ws.Range("I53").Value = 911.625
$ws.Range("J53").Value = 12643.111
$ws.Range("K53").Value = 911.625
$ws.Range("L53").Value = 12643.111
$ws.Range("M53").Value = -274.625
$ws.Range("N53").Value = -13917.111
$ws.Range("H113").Value = 14927.889
$ws.Range("I113").Value = 18558.834
$ws.Range("K113").Value = 18558.834
$ws.Range("M113").Value = -15304.834
$ws.Range("H125").Value = 15333
$ws.Range("J125").Value = 8000
$ws.Range("L125").Value = 72000
$ws.Range("N125").Value = -76920
$ws.Range("H138").Value = 1882.41
$ws.Range("I138").Value = 1271.5834
$ws.Range("J138").Value = 2075.3027
$ws.Range("K138").Value = 3814.7502
$ws.Range("L138").Value = 6225.908100000001
$ws.Range("M138").Value = 1325.2498
$ws.Range("N138").Value = -16505.9081
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1511502.9
$ws.Range("I19").Value = 1807801.6
$ws.Range("J19").Value = 30009
$ws.Range("K19").Value = 1807801.6
$ws.Range("L19").Value = 30009
$ws.Range("M19").Value = -1807572.6
$ws.Range("N19").Value = -30467
$ws.Range("H32").Value = 2294.6736
$ws.Range("I32").Value = 1768.5454
$ws.Range("K32").Value = 1768.5454
$ws.Range("M32").Value = -1481.5454
$ws.Range("H49").Value = 29326.666
$ws.Range("J49").Value = 29326.666
$ws.Range("L49").Value = 29326.666
$ws.Range("N49").Value = -29846.666
$ws.Range("H51").Value = 40000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H58").Value = 29000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H102").Value = 2701.7693
$ws.Range("I102").Value = 2760.25
$ws.Range("K102").Value = 2760.25
$ws.Range("M102").Value = -1138.25
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H139").Value = 75136.836
$ws.Range("J139").Value = 75136.836
$ws.Range("L139").Value = 75136.836
$ws.Range("N139").Value = -85416.836
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 28871.053
$ws.Range("I99").Value = 33223.848
$ws.Range("K99").Value = 33223.848
$ws.Range("M99").Value = -31725.848
$ws.Range("H107").Value = 2065.0588
$ws.Range("I107").Value = 2009.7273
$ws.Range("K107").Value = 2009.7273
$ws.Range("M107").Value = -89.72730000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3382.5454
$ws.Range("I132").Value = 3382.5454
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10147.6362
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7617.636200000001
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1450.91
$ws.Range("J131").Value = 1489.3043
$ws.Range("L131").Value = 4467.9129
$ws.Range("N131").Value = -14547.9129
$ws.Range("H132").Value = 1301.75
$ws.Range("J132").Value = 1489
$ws.Range("L132").Value = 13401
$ws.Range("N132").Value = -18461
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11998.577
$ws.Range("I80").Value = 10032.471
$ws.Range("K80").Value = 10032.471
$ws.Range("M80").Value = -9034.471
$ws.Range("H83").Value = 11998.577
$ws.Range("I83").Value = 10032.471
$ws.Range("K83").Value = 50162.355
$ws.Range("M83").Value = -45170.355
$ws.Range("H107").Value = 1822
$ws.Range("I107").Value = 1966.5
$ws.Range("K107").Value = 1966.5
$ws.Range("M107").Value = -46.5
$ws.Range("H123").Value = 75866.44
$ws.Range("J123").Value = 75866.44
$ws.Range("L123").Value = 75866.44
$ws.Range("N123").Value = -80766.44
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6989.5
$ws.Range("I40").Value = 2499.4
$ws.Range("K40").Value = 2499.4
$ws.Range("M40").Value = -2363.4
$ws.Range("H46").Value = 3506.6155
$ws.Range("I46").Value = 1865.6666
$ws.Range("K46").Value = 1865.6666
$ws.Range("M46").Value = -1677.6666
$ws.Range("H61").Value = 4068.7144
$ws.Range("I61").Value = 3079.05
$ws.Range("J61").Value = 6542.875
$ws.Range("K61").Value = 3079.05
$ws.Range("L61").Value = 6542.875
$ws.Range("M61").Value = -2877.05
$ws.Range("N61").Value = -6946.875
$ws.Range("H113").Value = 4068.7144
$ws.Range("I113").Value = 3079.05
$ws.Range("J113").Value = 6542.875
$ws.Range("K113").Value = 3079.05
$ws.Range("L113").Value = 6542.875
$ws.Range("M113").Value = -909.0500000000002
$ws.Range("N113").Value = -10882.875
$ws.Range("H122").Value = 6548.278
$ws.Range("I122").Value = 5949.9473
$ws.Range("J122").Value = 7217
$ws.Range("K122").Value = 17849.8419
$ws.Range("L122").Value = 21651
$ws.Range("M122").Value = -15399.8419
$ws.Range("N122").Value = -26551
$ws.Range("H132").Value = 14435.9375
$ws.Range("I132").Value = 1914.25
$ws.Range("J132").Value = 52001
$ws.Range("K132").Value = 5742.75
$ws.Range("L132").Value = 156003
$ws.Range("M132").Value = -3212.75
$ws.Range("N132").Value = -161063
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 15224
$ws.Range("I19").Value = 15000
$ws.Range("K19").Value = 15000
$ws.Range("M19").Value = -14826
$ws.Range("H96").Value = 1362.1333
$ws.Range("I96").Value = 1155.8
$ws.Range("J96").Value = 1465.3
$ws.Range("K96").Value = 1155.8
$ws.Range("L96").Value = 1465.3
$ws.Range("M96").Value = 217.2
$ws.Range("N96").Value = -4211.3
$ws.Range("H122").Value = 3705.2222
$ws.Range("I122").Value = 2817.745
$ws.Range("K122").Value = 8453.235000000001
$ws.Range("M122").Value = -6003.235000000001
$ws.Range("H132").Value = 7745.731
$ws.Range("I132").Value = 2318
$ws.Range("K132").Value = 6954
$ws.Range("M132").Value = -4424
